# Apply test data changes to avoid thread overwrite for test data values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new TestData/Browser values (TD6 / firefox) and Status -> Skip
$ws.Range("B3").Value = "TD6"
$ws.Range("C3").Value = "firefox"
$ws.Range("D3").Value = "Skip"

# Row 4: TestCase -> Test2, TestData/Browser -> Test2_TD1 / safari, Status stays Fail
$ws.Range("A4").Value = "Test2"
$ws.Range("B4").Value = "Test2_TD1"
$ws.Range("C4").Value = "safari"
$ws.Range("D4").Value = "Fail"
